$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.44
$ws.Range("D2").Value = 1.28

# Row 3
$ws.Range("D3").Value = 1.39
$ws.Range("E3").Value = 1.28
$ws.Range("G3").Value = 0.54

# Row 4
$ws.Range("B4").Value = 1.46
$ws.Range("C4").Value = 1.45
$ws.Range("E4").Value = 1.22
$ws.Range("F4").Value = 1

# Row 5
$ws.Range("C5").Value = 1.36
$ws.Range("D5").Value = 1.34
$ws.Range("G5").Value = 0.66

# Row 6
$ws.Range("D6").Value = 1.6

# Row 7
$ws.Range("C7").Value = 2.38
$ws.Range("E7").Value = 1.93
